$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.08258366666666667
$ws.Range("H2").Value = 0.247751
$ws.Range("M2").Value = 6.755097
$ws.Range("N2").Value = 20.265291
$ws.Range("O2").Value = 0.262181130417721
$ws.Range("P2").Value = 0.262181130417721
$ws.Range("Q2").Value = 0.557860678949
$ws.Range("R2").Value = 5.020746110541
$ws.Range("S2").Value = 0.262181130417721
$ws.Range("T2").Value = 0.262181130417721

# Row 3
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.08258366666666667
$ws.Range("H3").Value = 0.247751
$ws.Range("O3").Value = 0.179964029239562
$ws.Range("P3").Value = 0.179964029239562
$ws.Range("Q3").Value = 0.3829217433688889
$ws.Range("R3").Value = 3.446295690319999
$ws.Range("S3").Value = 0.179964029239562
$ws.Range("T3").Value = 0.179964029239562

# Row 4
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.08258366666666667
$ws.Range("H4").Value = 0.247751
$ws.Range("M4").Value = 12.28762933333333
$ws.Range("N4").Value = 36.862888
$ws.Range("O4").Value = 0.4769116637062769
$ws.Range("P4").Value = 0.4769116637062769
$ws.Range("Q4").Value = 1.014757484987556
$ws.Range("R4").Value = 9.132817364888
$ws.Range("S4").Value = 0.4769116637062769
$ws.Range("T4").Value = 0.4769116637062769

# Row 5
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.08258366666666667
$ws.Range("H5").Value = 0.247751
$ws.Range("M5").Value = 2.085501
$ws.Range("N5").Value = 6.256503
$ws.Range("O5").Value = 0.08094317663644024
$ws.Range("P5").Value = 0.08094317663644023
$ws.Range("Q5").Value = 0.172228319417
$ws.Range("R5").Value = 1.550054874753
$ws.Range("S5").Value = 0.08094317663644024
$ws.Range("T5").Value = 0.08094317663644023
